$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 10 de Septiembre de 2020 a las 03:18'
$ws.Range("B4").Value = 6549475
$ws.Range("C4").Value = 35244
$ws.Range("D4").Value = 3846095
$ws.Range("E4").Value = 2508141
$ws.Range("G4").Value = 1209
$ws.Range("H4").Value = 195239
$ws.Range("D5").Value = 3469084
$ws.Range("E5").Value = 918790
$ws.Range("B13").Value = 512293
$ws.Range("C13").Value = 12259
$ws.Range("E13").Value = 119145
$ws.Range("G13").Value = 253
$ws.Range("H13").Value = 10658
$ws.Range("B29").Value = 134294
$ws.Range("C29").Value = 546
$ws.Range("D29").Value = 118271
$ws.Range("E29").Value = 6868
$ws.Range("A55").Value = 'Venezuela'
$ws.Range("B55").Value = 56751
$ws.Range("C55").Value = 1188
$ws.Range("D55").Value = 45318
$ws.Range("E55").Value = 10981
$ws.Range("G55").Value = 8
$ws.Range("H55").Value = 452
$ws.Range("A56").Value = 'Nigeria'
$ws.Range("B56").Value = 55632
$ws.Range("C56").Value = 176
$ws.Range("D56").Value = 43610
$ws.Range("E56").Value = 10952
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 1070
$ws.Range("B76").Value = 25026
$ws.Range("C76").Value = 812
$ws.Range("D76").Value = 12309
$ws.Range("E76").Value = 12243
$ws.Range("G76").Value = 11
$ws.Range("H76").Value = 474
$ws.Range("E115").Value = 921
$ws.Range("H115").Value = 83
$ws.Range("A143").Value = 'Guadalupe'
$ws.Range("B143").Value = 2287
$ws.Range("C143").Value = 924
$ws.Range("D143").Value = 336
$ws.Range("E143").Value = 1928
$ws.Range("G143").Value = 5
$ws.Range("H143").Value = 23
$ws.Range("A144").Value = 'Guinea-Bisau'
$ws.Range("B144").Value = 2245
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 1127
$ws.Range("E144").Value = 1080
$ws.Range("H144").Value = 38
$ws.Range("A145").Value = 'Benin'
$ws.Range("B145").Value = 2242
$ws.Range("C145").Value = 29
$ws.Range("D145").Value = 1793
$ws.Range("E145").Value = 409
$ws.Range("H145").Value = 40
$ws.Range("A146").Value = 'Malta'
$ws.Range("B146").Value = 2162
$ws.Range("C146").Value = 63
$ws.Range("D146").Value = 1760
$ws.Range("E146").Value = 388
$ws.Range("H146").Value = 14
$ws.Range("A147").Value = 'Islandia'
$ws.Range("B147").Value = 2153
$ws.Range("C147").Value = 3
$ws.Range("D147").Value = 2067
$ws.Range("E147").Value = 76
$ws.Range("H147").Value = 10
$ws.Range("A148").Value = 'Botsuana'
$ws.Range("B148").Value = 2126
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 493
$ws.Range("E148").Value = 1624
$ws.Range("H148").Value = 9
$ws.Range("A149").Value = 'Sierra Leona'
$ws.Range("B149").Value = 2067
$ws.Range("C149").Value = 3
$ws.Range("D149").Value = 1622
$ws.Range("E149").Value = 373
$ws.Range("H149").Value = 72
$ws.Range("A150").Value = 'Yemen'
$ws.Range("B150").Value = 1999
$ws.Range("C150").Value = 5
$ws.Range("D150").Value = 1209
$ws.Range("E150").Value = 214
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 576
$ws.Range("A151").Value = 'Birmania'
$ws.Range("B151").Value = 1889
$ws.Range("C151").Value = 180
$ws.Range("D151").Value = 553
$ws.Range("E151").Value = 1324
$ws.Range("G151").Value = 2
$ws.Range("H151").Value = 12
$ws.Range("A152").Value = 'Nueva Zelanda'
$ws.Range("B152").Value = 1788
$ws.Range("C152").Value = 6
$ws.Range("D152").Value = 1639
$ws.Range("E152").Value = 125
$ws.Range("H152").Value = 24
$ws.Range("A153").Value = 'Georgia'
$ws.Range("B153").Value = 1773
$ws.Range("C153").Value = 44
$ws.Range("D153").Value = 1325
$ws.Range("E153").Value = 429
$ws.Range("H153").Value = 19
$ws.Range("A154").Value = 'Uruguay'
$ws.Range("B154").Value = 1741
$ws.Range("C154").Value = 29
$ws.Range("D154").Value = 1478
$ws.Range("E154").Value = 218
$ws.Range("H154").Value = 45
$ws.Range("A155").Value = 'Guyana'
$ws.Range("B155").Value = 1703
$ws.Range("C155").Value = 90
$ws.Range("D155").Value = 1045
$ws.Range("E155").Value = 610
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 48
$ws.Range("A156").Value = 'Togo'
$ws.Range("B156").Value = 1528
$ws.Range("C156").Value = 15
$ws.Range("D156").Value = 1144
$ws.Range("E156").Value = 348
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 36
$ws.Range("A157").Value = 'Republica de Chipre'
$ws.Range("B157").Value = 1514
$ws.Range("C157").Value = 3
$ws.Range("D157").Value = 1237
$ws.Range("E157").Value = 255
$ws.Range("H157").Value = 22
$ws.Range("A158").Value = 'Burkina Faso'
$ws.Range("B158").Value = 1476
$ws.Range("C158").Value = 10
$ws.Range("D158").Value = 1118
$ws.Range("E158").Value = 302
$ws.Range("H158").Value = 56
$ws.Range("A159").Value = 'Letonia'
$ws.Range("B159").Value = 1443
$ws.Range("C159").Value = 11
$ws.Range("D159").Value = 1234
$ws.Range("E159").Value = 174
$ws.Range("H159").Value = 35
$ws.Range("D168").Value = 862
$ws.Range("E168").Value = 21
$ws.Range("B170").Value = 722
$ws.Range("C170").Value = 6
$ws.Range("D170").Value = 662
$ws.Range("E170").Value = 18
$ws.Range("A214").Value = 'Montserrat'
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
$ws.Range("A215").Value = 'Islas Malvinas'
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
